# Added Test Data For Russia/Finland Market
#
# The workbook has one worksheet per market/country, each built from the
# same template (header rows, merged cells, shared-string lookups, etc.).
# The last sheet ("Turkey") is the template for two new market sheets:
# "Finland" and "Russia", appended at the end of the tab strip.

$wb = $excel.ActiveWorkbook

# --- Finland --------------------------------------------------------------
$turkey = $wb.Worksheets.Item("Turkey")
[void]$turkey.Copy($null, $turkey)

$finland = $wb.Worksheets.Item($wb.Worksheets.Count)
$finland.Name = "Finland"
$finland.Range("B2").Value = "Finland Market"
$finland.Range("B4").Value = "NGC-3130/T2949/T2888/T2945"
[void]$finland.Range("H12").Select()

# --- Russia -----------------------------------------------------------------
[void]$finland.Copy($null, $finland)

$russia = $wb.Worksheets.Item($wb.Worksheets.Count)
$russia.Name = "Russia"
$russia.Range("B2").Value = "Russia Market"
$russia.Range("B4").Value = "NGC-2929/T2917/T2911/T2902"
[void]$russia.Range("I12").Select()
